$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(8).Insert()
$ws.Range("H2").Value = "{spare:quantity}"
$ws.Range("H2").Font.Name = "Cambria"
Write-Output "done"
